$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2225774455180299
$ws.Range("C2").Value = 0.6025194276978968
$ws.Range("D2").Value = 0.9375366183193465
$ws.Range("E2").Value = 0.968264745986007
$ws.Range("F2").Value = 0.9517120452360237
$ws.Range("G2").Value = 51

# Row 3
$ws.Range("B3").Value = 0.5002927602626792
$ws.Range("C3").Value = 0.8093331159883599
$ws.Range("D3").Value = 2.512481039783275
$ws.Range("E3").Value = 1.585080767589865
$ws.Range("F3").Value = 1.519327252374612
$ws.Range("G3").Value = 50

# Row 4
$ws.Range("B4").Value = 0.7210151913624349
$ws.Range("C4").Value = 1.262988772388277
$ws.Range("D4").Value = 6.435165304488891
$ws.Range("E4").Value = 2.536762760781719
$ws.Range("F4").Value = 2.457343660326128
$ws.Range("G4").Value = 49

# Row 5
$ws.Range("B5").Value = 0.4586541891000874
$ws.Range("C5").Value = 1.384892259048143
$ws.Range("D5").Value = 7.994179379982225
$ws.Range("E5").Value = 2.82739798754654
$ws.Range("F5").Value = 2.819473144954236
$ws.Range("G5").Value = 48

# Row 6
$ws.Range("B6").Value = 0.4309535499285895
$ws.Range("C6").Value = 1.407201242206765
$ws.Range("D6").Value = 8.132525087741485
$ws.Range("E6").Value = 2.85175824496774
$ws.Range("F6").Value = 2.84948429316084
$ws.Range("G6").Value = 47

# Row 7
$ws.Range("B7").Value = 0.436267056338538
$ws.Range("C7").Value = 1.569053982142953
$ws.Range("D7").Value = 9.223627919858613
$ws.Range("E7").Value = 3.037042627270584
$ws.Range("F7").Value = 3.045889392437618
$ws.Range("G7").Value = 38

# Row 8
$ws.Range("B8").Value = 0.4110347441899493
$ws.Range("C8").Value = 1.588476331609555
$ws.Range("D8").Value = 9.491981118609186
$ws.Range("E8").Value = 3.080905892527259
$ws.Range("F8").Value = 3.095481328727775
$ws.Range("G8").Value = 37

# Row 9
$ws.Range("B9").Value = 0.08076273096934439
$ws.Range("C9").Value = 2.109961141030493
$ws.Range("D9").Value = 15.00836014869273
$ws.Range("E9").Value = 3.874062486420776
$ws.Range("F9").Value = 3.973840450533185
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = -0.6203614576892605
$ws.Range("C10").Value = 2.00672458876966
$ws.Range("D10").Value = 12.86991326340611
$ws.Range("E10").Value = 3.587466134112782
$ws.Range("F10").Value = 3.677701320433923
$ws.Range("G10").Value = 13

# Row 11 (G11 unchanged at 5)
$ws.Range("B11").Value = 1.509871661924674
$ws.Range("C11").Value = 1.524567153826987
$ws.Range("D11").Value = 3.867561475191954
$ws.Range("E11").Value = 1.966611673714959
$ws.Range("F11").Value = 1.408833311515586
